$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from refreshed cryptos feed.
# Numeric-looking text values (e.g. "590.24") must keep their
# original Text storage, so number format is pinned to "@"
# before assignment for those specific cells to avoid Excel
# auto-converting them into floating point numbers.

$ws.Range('D2').Value = '67.129.48'
$ws.Range('E2').Value = '  -5.05%  '
$ws.Range('D3').Value = '3.250.14'
$ws.Range('E3').Value = '  -7.85%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.24'
$ws.Range('E5').Value = '  -5.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.80'
$ws.Range('E6').Value = '  -12.55%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.241.31'
$ws.Range('E8').Value = '  -8.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  -11.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  -13.60%  '
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.507'
$ws.Range('E12').Value = '  -12.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.40'
$ws.Range('E13').Value = '  -17.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000243'
$ws.Range('E14').Value = '  -12.28%  '
$ws.Range('D15').Value = '3.771.08'
$ws.Range('E15').Value = '  -7.98%  '
$ws.Range('D16').Value = '67.204.51'
$ws.Range('E16').Value = '  -5.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '546.43'
$ws.Range('E17').Value = '  -10.28%  '
$ws.Range('D18').Value = '3.256.86'
$ws.Range('E18').Value = '  -7.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.29'
$ws.Range('E19').Value = '  -13.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.114'
$ws.Range('E20').Value = '  -6.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.22'
$ws.Range('E21').Value = '  -14.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.766'
$ws.Range('E22').Value = '  -13.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.83'
$ws.Range('E23').Value = '  -13.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.73'
$ws.Range('E24').Value = '  -13.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.56'
$ws.Range('E25').Value = '  -13.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.22'
$ws.Range('E27').Value = '  -14.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.09'
$ws.Range('E28').Value = '  -10.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '29.48'
$ws.Range('E29').Value = '  -12.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.14'
$ws.Range('E30').Value = '  -16.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('E31').Value = '  -11.15%  '
$ws.Range('E32').Value = '  -12.32%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.66'
$ws.Range('E33').Value = '  -17.86%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '544.16'
$ws.Range('E34').Value = '  -14.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.74'
$ws.Range('E35').Value = '  -15.47%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0445'
$ws.Range('E37').Value = '  -6.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.56'
$ws.Range('E38').Value = '  -5.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0857'
$ws.Range('E39').Value = '  -13.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.22'
$ws.Range('E40').Value = '  -14.56%  '
$ws.Range('D42').Value = '2.939.67'
$ws.Range('E42').Value = '  -12.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.63'
$ws.Range('E43').Value = '  -23.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.263'
$ws.Range('E44').Value = '  -15.41%  '
$ws.Range('D45').Value = '0.0₃0581'
$ws.Range('E45').Value = '  -19.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.55'
$ws.Range('E46').Value = '  -16.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.15'
$ws.Range('E47').Value = '  -15.20%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '127.51'
$ws.Range('E49').Value = '  -4.51%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.35'
$ws.Range('E50').Value = '  -20.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.114'
$ws.Range('E51').Value = '  -12.28%  '
